# IndividualWorkSheet.xlsx - "Add files via upload"
# Adds a new row (row 7) to the "이정원" worksheet describing the
# "Storage.java 코드 스멜 제거" task, mirroring the existing row 6
# ("Code Smell 제거" / Menu.java) entry, and moves the active selection
# from E6 to E7.

$wb = $excel.ActiveWorkbook

# The edited sheet is the 4th / last tab ("이정원"), which is also the
# workbook's active sheet.
$ws = $wb.ActiveSheet

# --- Row 7 values -----------------------------------------------------
# A7 reuses the same "Code Smell 제거" label already used by A6.
$ws.Range("A7").Value = "Code Smell 제거"
$ws.Range("B7").Value = "Storage.java 의 코드 스멜 1차 제거"
$ws.Range("C7").Value = 43603
$ws.Range("D7").Value = 43603
$ws.Range("E7").Value = "개선된 코드 Storage.java를 git commit함"

# --- Formatting ---------------------------------------------------------
# Match the wrap/vertical-top text formatting used by the surrounding
# rows (text columns wrap + align to top, date columns use the
# yyyy-mm-dd date format / top alignment).
$textCells = $ws.Range("A7,B7,E7")
$textCells.WrapText = $true
$textCells.VerticalAlignment = -4160   # xlTop

$dateCells = $ws.Range("C7,D7")
$dateCells.NumberFormat = "yyyy\-mm\-dd;@"
$dateCells.VerticalAlignment = -4160   # xlTop

# Row grows to fit the now-wrapped text (matches rows 2-4/6).
$ws.Rows.Item(7).RowHeight = 35

# --- Selection ----------------------------------------------------------
# The saved selection moves from E6 to the newly-filled E7.
$ws.Range("E7").Select()
